# Insert a new weekly price record for Albahaca right before the current
# row 108, pushing the existing rows 108-192 down to 109-193.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 108 (existing row 108 and below shift to 109+)
$ws.Rows.Item(108).Insert()

# The data that used to live in row 108 is now in row 109; use it as the
# template for the new row's constant columns (Mercado, Region, etc.)
$ws.Cells.Item(108, 1).Value = $ws.Cells.Item(109, 1).Value()
$ws.Cells.Item(108, 2).Value = $ws.Cells.Item(109, 2).Value()
$ws.Cells.Item(108, 3).Value = $ws.Cells.Item(109, 3).Value()
$ws.Cells.Item(108, 4).Value = 45096
$ws.Cells.Item(108, 5).Value = $ws.Cells.Item(109, 5).Value()
$ws.Cells.Item(108, 6).Value = $ws.Cells.Item(109, 6).Value()
$ws.Cells.Item(108, 7).Value = $ws.Cells.Item(109, 7).Value()
$ws.Cells.Item(108, 8).Value = $ws.Cells.Item(109, 8).Value()
$ws.Cells.Item(108, 9).Value = $ws.Cells.Item(109, 9).Value()
$ws.Cells.Item(108, 10).Value = 600
$ws.Cells.Item(108, 11).Value = 4000
$ws.Cells.Item(108, 12).Value = 5000
$ws.Cells.Item(108, 13).Value = 4500
$ws.Cells.Item(108, 14).Value = $ws.Cells.Item(109, 14).Value()
$ws.Cells.Item(108, 15).Value = $ws.Cells.Item(109, 15).Value()
$ws.Cells.Item(108, 16).Value = 4500
$ws.Cells.Item(108, 17).Value = $ws.Cells.Item(109, 17).Value()
$ws.Cells.Item(108, 18).Value = $ws.Cells.Item(109, 18).Value()
